# Update the cached "datetimeFigureOut" field text on every Date
# Placeholder shape across all slide masters and all of their slide
# layouts, from "21-11-2022" to "22-11-2022".

$p = $ppt.ActivePresentation

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "21-11-2022") {
                $tr.Text = "22-11-2022"
            }
        }
    }
}

for ($d = 1; $d -le $p.Designs.Count; $d++) {
    $design = $p.Designs.Item($d)
    $master = $design.SlideMaster

    Update-DateShapes $master.Shapes

    for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
        $layout = $master.CustomLayouts.Item($l)
        Update-DateShapes $layout.Shapes
    }
}

Write-Host "Updated cached date field text on all masters/layouts."
